$d = $word.ActiveDocument

# 1) "App Name" -> "NUSTeats" (keep the surrounding smart quotes), and make Word
#    split the text into 3 runs just like a real retype-over-selection would:
#    quote | NUSTeats | quote
$r = $d.Content
[void]$r.Find.Execute("App Name")
$start = $r.Start
$r.Text = "NUSTeats"
$mid = $d.Range($start, $start + 8)
$mid.Font.Bold = 1
$mid.Font.Bold = 0

# 2) Remove the "Time Filter" bullet entirely (list collapses by one item).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Time Filter`r") {
        $p.Range.Delete()
        break
    }
}

# 3) "UI: Qt" -> "UI: Raylib" (split the replaced word into its own run).
$r2 = $d.Content
[void]$r2.Find.Execute(" Qt")
$spaceStart = $r2.Start
$wordRange = $d.Range($spaceStart + 1, $r2.End)
$wordRange.Text = "Raylib"
$mid2 = $d.Range($spaceStart + 1, $spaceStart + 1 + 6)
$mid2.Font.Bold = 1
$mid2.Font.Bold = 0

# 4) "mySQL" -> "mySQL(If data is too much to handle in 2D arrays)"
$r3 = $d.Content
[void]$r3.Find.Execute("mySQL")
$r3.Text = "mySQL(If data is too much to handle in 2D arrays)"
